$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 442.7143
$ws.Range("I33").Value = 442.7143
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 442.7143
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -213.7143
$ws.Range("H57").Value = 64044.332
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 64044.332
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 192132.996
$ws.Range("N57").Value = -193130.996
$ws.Range("H74").Value = 11707.56
$ws.Range("I74").Value = 11921.261
$ws.Range("J74").Value = 9250
$ws.Range("K74").Value = 11921.261
$ws.Range("L74").Value = 9250
$ws.Range("M74").Value = -10985.261
$ws.Range("H77").Value = 11707.56
$ws.Range("I77").Value = 11921.261
$ws.Range("J77").Value = 9250
$ws.Range("K77").Value = 59606.305
$ws.Range("L77").Value = 46250
$ws.Range("M77").Value = -54926.305
$ws.Range("H111").Value = 4134.9
$ws.Range("I111").Value = 5335.7144
$ws.Range("J111").Value = 1333
$ws.Range("K111").Value = 16007.1432
$ws.Range("L111").Value = 3999
$ws.Range("M111").Value = -12940.1432
$ws.Range("H115").Value = 376.4
$ws.Range("I115").Value = 370.5
$ws.Range("J115").Value = 400
$ws.Range("K115").Value = 1111.5
$ws.Range("L115").Value = 1200
$ws.Range("M115").Value = 455.5
$ws.Range("N115").Value = -4334
$ws.Range("H116").Value = 7500.8335
$ws.Range("I116").Value = 4335
$ws.Range("J116").Value = 10666.667
$ws.Range("K116").Value = 4335
$ws.Range("L116").Value = 10666.667
$ws.Range("M116").Value = -893
$ws.Range("N116").Value = -17550.667
$ws.Range("H132").Value = 7401.4736
$ws.Range("I132").Value = 7625.1763
$ws.Range("J132").Value = 5500
$ws.Range("K132").Value = 22875.5289
$ws.Range("L132").Value = 16500
$ws.Range("M132").Value = -20345.5289

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2357.6
$ws.Range("I2").Value = 2299.8462
$ws.Range("J2").Value = 2464.8572
$ws.Range("K2").Value = 2299.8462
$ws.Range("L2").Value = 2464.8572
$ws.Range("M2").Value = -2186.8462
$ws.Range("H32").Value = 1883107
$ws.Range("I32").Value = 843653.8
$ws.Range("J32").Value = 47619050
$ws.Range("K32").Value = 843653.8
$ws.Range("L32").Value = 47619050
$ws.Range("M32").Value = -843366.8
$ws.Range("H116").Value = 2357.6
$ws.Range("I116").Value = 2299.8462
$ws.Range("J116").Value = 2464.8572
$ws.Range("K116").Value = 2299.8462
$ws.Range("L116").Value = 2464.8572
$ws.Range("M116").Value = -5.846199999999953
$ws.Range("H132").Value = 2286
$ws.Range("I132").Value = 1186.1
$ws.Range("J132").Value = 5428.5713
$ws.Range("K132").Value = 3558.3
$ws.Range("L132").Value = 16285.7139
$ws.Range("M132").Value = -1028.3

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2357.6
$ws.Range("I3").Value = 2299.8462
$ws.Range("J3").Value = 2464.8572
$ws.Range("K3").Value = 2299.8462
$ws.Range("L3").Value = 2464.8572
$ws.Range("M3").Value = -2185.8462
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H86").Value = 2661.4194
$ws.Range("I86").Value = 2608.5
$ws.Range("J86").Value = 2842.8572
$ws.Range("K86").Value = 2608.5
$ws.Range("L86").Value = 2842.8572
$ws.Range("M86").Value = -1485.5
$ws.Range("H89").Value = 2661.4194
$ws.Range("I89").Value = 2608.5
$ws.Range("J89").Value = 2842.8572
$ws.Range("K89").Value = 13042.5
$ws.Range("L89").Value = 14214.286
$ws.Range("M89").Value = -7426.5
$ws.Range("H105").Value = 8967231
$ws.Range("I105").Value = 477671.62
$ws.Range("J105").Value = 31252324
$ws.Range("K105").Value = 477671.62
$ws.Range("L105").Value = 31252324
$ws.Range("M105").Value = -475924.62
$ws.Range("N105").Value = -31255818

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5452.231
$ws.Range("I86").Value = 4796.5
$ws.Range("J86").Value = 7638
$ws.Range("K86").Value = 4796.5
$ws.Range("L86").Value = 7638
$ws.Range("M86").Value = -3673.5
$ws.Range("N86").Value = -9884
$ws.Range("H89").Value = 5452.231
$ws.Range("I89").Value = 4796.5
$ws.Range("J89").Value = 7638
$ws.Range("K89").Value = 23982.5
$ws.Range("L89").Value = 38190
$ws.Range("M89").Value = -18366.5
$ws.Range("N89").Value = -49422
$ws.Range("H134").Value = 3753.1667
$ws.Range("I134").Value = 3959.35
$ws.Range("J134").Value = 3340.8
$ws.Range("K134").Value = 11878.05
$ws.Range("L134").Value = 10022.4
$ws.Range("M134").Value = -9343.049999999999
$ws.Range("H140").Value = 69962.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 69962.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 69962.5
$ws.Range("N140").Value = -80322.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 87.22
$ws.Range("I4").Value = 68.15625
$ws.Range("J4").Value = 544.75
$ws.Range("K4").Value = 204.46875
$ws.Range("L4").Value = 1634.25
$ws.Range("M4").Value = -92.46875
$ws.Range("N4").Value = -1858.25
$ws.Range("H5").Value = 1291.4286
$ws.Range("I5").Value = 1265.25
$ws.Range("J5").Value = 1326.3334
$ws.Range("K5").Value = 3795.75
$ws.Range("L5").Value = 3979.0002
$ws.Range("M5").Value = -3683.75
$ws.Range("N5").Value = -4203.0002
$ws.Range("H113").Value = 1033.875
$ws.Range("I113").Value = 874.8
$ws.Range("J113").Value = 1299
$ws.Range("K113").Value = 2624.4
$ws.Range("L113").Value = 3897
$ws.Range("M113").Value = -454.3999999999996
$ws.Range("N113").Value = -8237
$ws.Range("H114").Value = 2745.2144
$ws.Range("I114").Value = 1941.5
$ws.Range("J114").Value = 3066.7
$ws.Range("K114").Value = 5824.5
$ws.Range("L114").Value = 9200.099999999999
$ws.Range("M114").Value = -2570.5
$ws.Range("N114").Value = -15708.1
$ws.Range("H116").Value = 4000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 12000
$ws.Range("N116").Value = -18884
$ws.Range("M116").ClearContents()
$ws.Range("H117").Value = 1320.75
$ws.Range("I117").Value = 561
$ws.Range("J117").Value = 3600
$ws.Range("K117").Value = 1683
$ws.Range("L117").Value = 10800
$ws.Range("M117").Value = 1759
$ws.Range("N117").Value = -17684
$ws.Range("H131").Value = 3278530
$ws.Range("I131").Value = 27338.334
$ws.Range("J131").Value = 4904126
$ws.Range("K131").Value = 82015.00199999999
$ws.Range("L131").Value = 14712378
$ws.Range("M131").Value = -76975.00199999999
$ws.Range("N131").Value = -14722458
$ws.Range("H132").Value = 1138.5
$ws.Range("I132").Value = 1150.75
$ws.Range("J132").Value = 1114
$ws.Range("K132").Value = 10356.75
$ws.Range("L132").Value = 10026
$ws.Range("M132").Value = -7826.75
$ws.Range("H134").Value = 3221.3333
$ws.Range("H135").Value = 1291.4286
$ws.Range("I135").Value = 1265.25
$ws.Range("J135").Value = 1326.3334
$ws.Range("K135").Value = 11387.25
$ws.Range("L135").Value = 11937.0006
$ws.Range("M135").Value = -8852.25
$ws.Range("N135").Value = -17007.0006
$ws.Range("H136").Value = 6664.75
$ws.Range("I136").Value = 2220
$ws.Range("J136").Value = 19999
$ws.Range("K136").Value = 6660
$ws.Range("L136").Value = 59997
$ws.Range("M136").Value = -1560
$ws.Range("N136").Value = -70197
$ws.Range("H137").Value = 2322.6
$ws.Range("I137").Value = 2322.6
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6967.799999999999
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1867.799999999999
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 7204.4165
$ws.Range("I139").Value = 7772.6113
$ws.Range("J139").Value = 5499.8335
$ws.Range("K139").Value = 23317.8339
$ws.Range("L139").Value = 16499.5005
$ws.Range("M139").Value = -18177.8339
$ws.Range("H140").Value = 2226.7368
$ws.Range("I140").Value = 2083.7778
$ws.Range("J140").Value = 4800
$ws.Range("K140").Value = 6251.3334
$ws.Range("L140").Value = 14400
$ws.Range("M140").Value = -1071.3334
$ws.Range("H141").Value = 19999
$ws.Range("I141").Value = 19999
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 59997
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -54817
$ws.Range("N141").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 39349.4
$ws.Range("I57").Value = 3500
$ws.Range("J57").Value = 43332.668
$ws.Range("K57").Value = 3500
$ws.Range("L57").Value = 43332.668
$ws.Range("M57").Value = -2680
$ws.Range("N57").Value = -44972.668
$ws.Range("H122").Value = 37138176
$ws.Range("I122").Value = 4275683.5
$ws.Range("J122").Value = 90913170
$ws.Range("K122").Value = 12827050.5
$ws.Range("L122").Value = 272739510
$ws.Range("M122").Value = -12824600.5
$ws.Range("H132").Value = 1734.2354
$ws.Range("I132").Value = 1648.875
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 4946.625
$ws.Range("L132").Value = 9300
$ws.Range("M132").Value = -2416.625

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8015.643
$ws.Range("I136").Value = 3521.4
$ws.Range("J136").Value = 19251.25
$ws.Range("K136").Value = 10564.2
$ws.Range("L136").Value = 57753.75
$ws.Range("M136").Value = -8014.200000000001
$ws.Range("H139").Value = 73131
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 73131
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 73131
$ws.Range("N139").Value = -83411
$ws.Range("H140").Value = 101722.73
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 101722.73
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 101722.73
$ws.Range("N140").Value = -112082.73

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 166667970
$ws.Range("I100").Value = 1578.75
$ws.Range("J100").Value = 500000740
$ws.Range("K100").Value = 3157.5
$ws.Range("L100").Value = 1000001480
$ws.Range("M100").Value = -2616.5
$ws.Range("N100").Value = -1000002562
$ws.Range("H136").Value = 9806530
$ws.Range("I136").Value = 10419251
$ws.Range("J136").Value = 3004.5
$ws.Range("K136").Value = 31257753
$ws.Range("L136").Value = 9013.5
$ws.Range("M136").Value = -31255203
$ws.Range("N136").Value = -14113.5
$ws.Range("H139").Value = 79998.73
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 79998.73
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 79998.73
$ws.Range("N139").Value = -90278.73
